$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-14 (B column text updates due to list reorder/rename/removal)
$ws.Range("B2").Value = "صيدلية النورس"
$ws.Range("B3").Value = "صيدلية د.أحمد فتح الله"
$ws.Range("B4").Value = "صيدلية د.حنان يحيى"
$ws.Range("B5").Value = "صيدلية د.سامي سمير"
$ws.Range("B6").Value = "صيدلية د.شريف نوح"
$ws.Range("B7").Value = "صيدلية د.كامل مرعي"
$ws.Range("B8").Value = "صيدلية د.محمد فتح الله"
$ws.Range("B9").Value = "صيدلية د.ناصر الشافعي"
$ws.Range("B10").Value = "صيدلية د.هاني شاكر"
$ws.Range("B11").Value = "صيدلية د.هاني عبد السلام"
$ws.Range("B12").Value = "صيدلية د.وفاء عبده"
$ws.Range("B13").Value = "صيدلية العقباوى القاهرة"
$ws.Range("B14").Value = "صيدلية د.وليد"

# Extend formatting (style) for new rows 15-29 by copying from row 14
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B29").PasteSpecial(-4122)

# Populate new rows 15-29
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "صيدلية د. عبد الله كامل - القاهرة"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "صيدلية د.هبة - القاهرة"
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "صيدلية د.أسماء - القاهرة"
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "صيدلية 70 فدان"
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "صيدلية العبور"
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "صيدلية د.احمد احمد"
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "صيدلية العبور - القاهرة"
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "صيدلية د.ابراهيم"
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "صيدلية الصفا والمروة"
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "صيدلية د. نعيم"
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "صيدلية د.اسماء - القاهرة"
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "صيدلية د.شيماء"
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "صيدلية /د.عادل سعيد"
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "كوزى كورنر"
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "صيدلية الحرية - القاهرة"

$excel.CutCopyMode = 0
Write-Host "done"
